$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 / Spring 2022 block (rows 4-10) ---
# Column A: Fall 2022 courses (PSYC 1101 removed, CPSC 3415 added)
$ws.Range("A4").Value = "POLS 1101"
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = "PSYC 1105"
$ws.Range("B5").Value = 2
$ws.Range("A6").Value = "DSCI 3111"
$ws.Range("B6").Value = 3
$ws.Range("A7").Value = "CPSC 3121"
$ws.Range("B7").Value = 3
$ws.Range("A8").Value = "CPSC 3165"
$ws.Range("B8").Value = 3
$ws.Range("A9").Value = "CPSC 3415"
$ws.Range("B9").Value = 1
$ws.Range("A10").Value = "CPSC 4000"
$ws.Range("B10").Value = 0

# Column C: Spring 2022 courses (GEOL 3225 and CYBR 4125 added at top)
$ws.Range("C4").Value = "GEOL 3225"
$ws.Range("D4").Value = 3
$ws.Range("C5").Value = "CYBR 4125"
$ws.Range("D5").Value = 3
$ws.Range("C6").Value = "CPSC 4135"
$ws.Range("D6").Value = 3
$ws.Range("C7").Value = "CPSC 4148"
$ws.Range("D7").Value = 3
$ws.Range("C8").Value = "CPSC 4155"
$ws.Range("D8").Value = 3

# --- Fall 2023 / Spring 2023 block (rows 13-14) ---
$ws.Range("A13").Value = "CPSC 4157"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3
$ws.Range("A14").Value = "CPSC 4175"
$ws.Range("B14").Value = 3
